# Finish tod adult norms, rerun grade for ORF in school-age sample
#
# Each of the 6 age-band lookup tabs (raw -> ss) gains one more row at the
# top of the table: "raw" now starts at 0 instead of 1, every existing raw
# score shifts down one position, and the ss column is rescored accordingly.
# The table grows from 20 data rows (A2:B21) to 21 data rows (A2:B22).

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "18.0-23.11" ---
$ws = $wb.Worksheets.Item(1)

$raw = @(0, 1, 2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20)
$ss  = @(40, 40, 40, 54, 67, 77, 85, 92, 98, 104, 109, 114, 119, 123, 127, 130, 130, 130, 130, 130, 130)

for ($i = 0; $i -lt $raw.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $raw[$i]
    $ws.Cells.Item($row, 2).Value = $ss[$i]
}

# --- Sheet 2: "24.0-39.11" ---
$ws = $wb.Worksheets.Item(2)

$raw = @(0, 1, 2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20)
$ss  = @(40, 40, 40, 54, 67, 77, 86, 93, 99, 105, 111, 116, 121, 126, 130, 130, 130, 130, 130, 130, 130)

for ($i = 0; $i -lt $raw.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $raw[$i]
    $ws.Cells.Item($row, 2).Value = $ss[$i]
}

# --- Sheet 3: "40.0-49.11" ---
$ws = $wb.Worksheets.Item(3)

$raw = @(0, 1, 2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20)
$ss  = @(40, 40, 40, 54, 67, 77, 85, 93, 99, 105, 111, 117, 122, 127, 130, 130, 130, 130, 130, 130, 130)

for ($i = 0; $i -lt $raw.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $raw[$i]
    $ws.Cells.Item($row, 2).Value = $ss[$i]
}

# --- Sheet 4: "50.0-59.11" ---
$ws = $wb.Worksheets.Item(4)

$raw = @(0, 1, 2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20)
$ss  = @(40, 40, 40, 56, 68, 78, 86, 93, 99, 105, 111, 116, 122, 127, 130, 130, 130, 130, 130, 130, 130)

for ($i = 0; $i -lt $raw.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $raw[$i]
    $ws.Cells.Item($row, 2).Value = $ss[$i]
}

# --- Sheet 5: "60.0-69.11" ---
$ws = $wb.Worksheets.Item(5)

$raw = @(0, 1, 2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20)
$ss  = @(40, 40, 40, 58, 70, 79, 87, 93, 100, 105, 111, 116, 121, 127, 130, 130, 130, 130, 130, 130, 130)

for ($i = 0; $i -lt $raw.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $raw[$i]
    $ws.Cells.Item($row, 2).Value = $ss[$i]
}

# --- Sheet 6: "70.0-89.11" ---
$ws = $wb.Worksheets.Item(6)

$raw = @(0, 1, 2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20)
$ss  = @(40, 40, 48, 65, 75, 83, 90, 96, 102, 107, 113, 118, 122, 127, 130, 130, 130, 130, 130, 130, 130)

for ($i = 0; $i -lt $raw.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $raw[$i]
    $ws.Cells.Item($row, 2).Value = $ss[$i]
}

